# 03_VPC.pptx edit
# 1. Retitle the chapter from "AWS Introduction" to "VPC" on the title slide
#    and every section-header slide that repeats the running chapter title.
# 2. Remove the old "2.1 What is Cloud?" slide (it belonged to the previous
#    chapter and no longer belongs in this deck).

$p = $ppt.ActivePresentation

# --- 1. Update the chapter title text on slide 1 ("3 AWS Introduction" -> "3 VPC") ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "3 VPC"

# --- 2. Update the repeated chapter-title banner on slides 2-10
#         ("2 AWS Introduction" -> "3 VPC") ---
for ($i = 2; $i -le 10; $i++) {
    $s = $p.Slides.Item($i)
    $s.Shapes.Item(1).TextFrame.TextRange.Text = "3 VPC"
}

# --- 3. Delete the stray "2.1 What is Cloud?" slide (position 11) ---
$oldSlide = $p.Slides.Item(11)
$oldSlide.Delete()

# --- 4. The slide that used to be #12 ("End of Chapter") is now #11; refresh
#         its cached slide-number field so it reads "11" instead of "12". ---
$lastSlide = $p.Slides.Item($p.Slides.Count)
for ($j = 1; $j -le $lastSlide.Shapes.Count; $j++) {
    $shp = $lastSlide.Shapes.Item($j)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "12") {
            $shp.TextFrame.TextRange.Text = "11"
        }
    }
}
